$wb = $excel.ActiveWorkbook

# Add a new worksheet named "REPORT" and move it to be the first sheet in the workbook.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "REPORT"
$newSheet.Move($wb.Worksheets.Item(1))

# Write the report text into cell A26 and select it (as it appears in the target file).
$newSheet.Range("A26").Value = "3. For the days selected in step two, identify all emails which denote a change detected. Open each email and ensure the details within (e.g. each change is captured within the monthly tracker (findings tab)."
$newSheet.Range("A26:I26").Select()
